$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: clear B2, D2, E2; set C2 to new value
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -0.5911246880189821
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 value tweaks
$ws.Range("B3").Value = -0.75226107008933984
$ws.Range("C3").Value = 0.81958017826096596
$ws.Range("D3").Value = -0.096784206100009193
$ws.Range("E3").Value = 2.3909692343347553

# Update selection to match new sqref
$ws.Range("B1:E3").Select()
